# Update "想去人数" (want-to-go count) figures for the two exhibition
# entries that appear on both the "展览" sheet and the "全部类型" sheet.
#   F2: 1583 -> 1585
#   F4: 40   -> 41

$wb = $excel.ActiveWorkbook
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1585
    $ws.Range("F4").Value = 41
}
